# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview + per-locale sheets: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn sheet: records the target + handback file links/dates
#  - de-de sheet: records the target + handback file links/dates
#  - widens a few columns that now hold longer content

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9da188b4639289da48b408f65a0dcbe456afe9d4/e2e/e99a8093-d901-4ade-b759-d1188cbcd08a.md"
$targetName = "e99a8093-d901-4ade-b759-d1188cbcd08a.md"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) on both data rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (I) now links back to the source markdown file
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetUrl, "", "", $targetName)
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetUrl, "", "", $targetName)
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = 15570276

# Latest Handback File (J) now references the generated xliff
$wsZh.Range("J2").Value = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.zh-cn.xlf"
$wsZh.Range("J3").Value = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZh.Range("K2").Value = "2016-08-29 03:02:41"
$wsZh.Range("K3").Value = "2016-08-29 03:02:41"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZh.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetUrl, "", "", $targetName)
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetUrl, "", "", $targetName)
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = 15570276

$wsDe.Range("J2").Value = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.de-de.xlf"
$wsDe.Range("J3").Value = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-29 03:02:48"
$wsDe.Range("K3").Value = "2016-08-29 03:02:48"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.16666666666667
